# The commit adds a new data row to the Users sheet: the user typed
# "asd" into cell A2 (directly below the existing header row) and then
# moved on to A3, which is exactly what the target worksheet reflects
# (A2 now holds the text "asd", and the sheet's active selection ends
# up on A3).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "asd"

# Mirror the natural post-edit cursor position (pressing Enter after
# typing into A2 moves the selection down to A3).
$ws.Range("A3").Select()
